$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = 45876.87520870371

$ws.Range("A24").NumberFormat = $ws.Range("A23").NumberFormat
$ws.Range("A24").Value = 45876.91711516142
$ws.Range("B24").Value = 2025
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 14.5
$ws.Range("E24").Value = 90.59999999999999
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 7.48
$ws.Range("H24").Value = "ESE"
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = "22:00:38"
